# fix NPC HP error
# Insert a new "MAXHP" column before column G, shifting existing columns
# G..Y to H..Z, and populate the new column with the correct HP values
# (same as the SalePrice/F column, which is what the HP values should
# have mirrored).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column G (shifts G:Y -> H:Z)
$ws.Columns("G:G").Insert()

# New column header
$ws.Range("G1").Value = "MAXHP"

# New column values - mirror the SalePrice column (fixes the NPC HP bug)
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("G6").Value = $ws.Range("F6").Value2

# Restore view state: scroll back to show column A and select I8
$ws.Range("I8").Select()
